$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# -----------------------------------------------------------------
# 1) "Scout Positions" table (Table 3): the Wise Guys Patrol Leader
#    row gets "(Sr.)" inserted into the position text.
# -----------------------------------------------------------------
$scoutTbl = $s.Shapes.Item("Table 3").Table
$wiseGuysCell = $scoutTbl.Cell(11, 2).Shape.TextFrame.TextRange
$wiseGuysCell.Text = "Wise Guys (Sr.)  Patrol Leader"

# -----------------------------------------------------------------
# 2) Committee table (Table 128): add a new "Paul Besser - Eagle
#    Committee Coordinator" row right after the Trek Coordinator row.
# -----------------------------------------------------------------
$commTbl = $s.Shapes.Item("Table 128").Table
$newRow = $commTbl.Rows.Add(11)
$commTbl.Cell(11, 1).Shape.TextFrame.TextRange.Text = "Paul Besser"
$commTbl.Cell(11, 2).Shape.TextFrame.TextRange.Text = "Eagle Committee Coordinator"
